$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price(D) and Volume(E) columns to text so numeric-looking values
# (prices, percentages) are not auto-converted to numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '49.345.31'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '2.624.07'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '112.11'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").Value = '323.41'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("D10").Value = '39.84'
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").Value = '19.81'
$ws.Range("E11").Value = '  -3.85%  '
$ws.Range("D12").Value = '0.0811'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '7.28'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '3.028.87'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '2.630.02'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '49.232.04'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").Value = '12.96'
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '269.83'
$ws.Range("E23").Value = '  -3.15%  '
$ws.Range("D24").Value = '68.59'
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").Value = '26.20'
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '4.08'
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '10.32'
$ws.Range("E29").Value = '  +3.70%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").Value = '  -3.51%  '
$ws.Range("D32").Value = '34.97'
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").Value = '5.48'
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("D35").Value = '0.0815'
$ws.Range("E35").Value = '  +3.16%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '19.01'
$ws.Range("E37").Value = '  -3.48%  '
$ws.Range("D38").Value = '4.91'
$ws.Range("E38").Value = '  +3.32%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").Value = '127.37'
$ws.Range("E41").Value = '  +3.15%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '22.36'
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.111'
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("D44").Value = '0.0320'
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("D46").Value = '2.060.82'
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").Value = '2.16'
$ws.Range("E47").Value = '  +7.59%  '
$ws.Range("D48").Value = '3.22'
$ws.Range("E48").Value = '  -3.99%  '
$ws.Range("D49").Value = '2.13'
$ws.Range("E49").Value = '  -8.49%  '
$ws.Range("D50").Value = '8.92'
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").Value = '58.99'
$ws.Range("E51").Value = '  +1.64%  '

# Restore default "Normal" style so no stray number-format style is left
# attached to the cells (matches original unstyled data cells).
$ws.Range("D2:E51").Style = "Normal"
